# feat: add 2022-Q4 data
#
# Before:  sheets = [ 总计 , 2022-Q2 ]
# After:   sheets = [ 总计 , 2022-Q4 , 2022-Q2 ]
#   - "2022-Q4" carries the fresh Q4 fund-holdings detail (reuses the
#     physical sheet that used to be named "2022-Q2")
#   - a brand-new "2022-Q2" tab is appended right after it, holding an
#     exact copy of the original Q2 fund-holdings detail
#   - the "总计" (totals) sheet gains a new row for 2022-Q4 and keeps the
#     old 2022-Q2 totals row (shifted down one row)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the existing "2022-Q2" sheet so we end up with two tabs:
#    the original (which we will turn into "2022-Q4") and a verbatim
#    copy (which keeps the "2022-Q2" data/name).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $q2)

# Rename the original sheet out of the way *before* renaming the copy,
# otherwise the copy's rename collides with the still-existing original.
$q2.Name = "2022-Q4"

$q2Copy = $wb.Worksheets.Item("2022-Q2 (2)")
$q2Copy.Name = "2022-Q2"

$q4 = $q2

# ---------------------------------------------------------------------
# 2. Replace the content of the (renamed) "2022-Q4" sheet with the new
#    quarter's fund-holdings data.
# ---------------------------------------------------------------------
$q4.Cells.Clear()

$zj = $wb.Worksheets.Item("总计")

# Match the header / first-column formatting used elsewhere in the
# workbook (bold, centred, thin border) by copying it across instead of
# re-describing it by hand.
$zj.Range("B1:D1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)
$zj.Range("A2").Copy()
$q4.Range("A2").PasteSpecial(-4122)

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4.Range("A2").Value = 0

# These columns hold numeric-looking text in the source data (fund code
# with a leading zero, percentages/ratios kept as plain strings) -
# force a text format first so Excel doesn't silently coerce them to
# numbers.
$q4.Range("B2:G2").NumberFormat = "@"
$q4.Range("B2").Value = "001797"
$q4.Range("C2").Value = "国新国证新利灵活配置混合"
$q4.Range("D2").Value = "0.02"
$q4.Range("E2").Value = "81.37"
$q4.Range("F2").Value = "3.91"
$q4.Range("G2").Value = "0.0008"
$q4.Range("H2").Value = 5

# Match the page margins used across the rest of the workbook.
$q4.PageSetup.LeftMargin = 54
$q4.PageSetup.RightMargin = 54
$q4.PageSetup.TopMargin = 72
$q4.PageSetup.BottomMargin = 72
$q4.PageSetup.HeaderMargin = 36
$q4.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 3. Update the "总计" (totals) sheet: push the existing 2022-Q2 totals
#    row down to row 3, then write the new 2022-Q4 totals into row 2.
# ---------------------------------------------------------------------
$zj.Range("A2:D2").Copy()
$zj.Range("A3:D3").PasteSpecial(-4122)

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q2"
$zj.Range("C3").Value = 3
$zj.Range("D3").Value = 0

$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0

$zj.Activate()
